# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E16:E23) is reversed into descending order
# (2106 down to 2011, replacing the previous ascending 2011 up to 2106),
# and the "Valor Mora" amounts in F16 / F23 swap (120000 <-> 84000) so the
# smallest value now sits on the new top row (2106) instead of the old
# bottom row (2106 was previously last / 2011 first).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2106"
$ws.Range("E17").Value = "2105"
$ws.Range("E18").Value = "2104"
$ws.Range("E19").Value = "2103"
$ws.Range("E20").Value = "2102"
$ws.Range("E21").Value = "2101"
$ws.Range("E22").Value = "2012"
$ws.Range("E23").Value = "2011"

$ws.Range("F16").Value = 84000
$ws.Range("F23").Value = 120000
